# "Generate Report for handoff"
# b.md.md has been handed off again: its status flips from the stale
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# per-locale sheets get a fresh "Latest Handoff File" / "Latest Handoff
# Datetime" pointing at the new xlf for b.md.md (hash
# b3a40d6229ff1a8b48804fcfc66c95922eb78fd0).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status column for b.md.md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- Per-locale detail sheets ---
$locales = @(
    @{ Sheet = "zh-cn"; File = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"; Datetime = "2016-02-16 14:50:38" },
    @{ Sheet = "de-de"; File = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"; Datetime = "2016-02-16 14:50:52" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Status: "Handed back: in sync with en-US" -> "Ready for handoff"
    $ws.Range("B3").Value = "Ready for handoff"

    # Latest Handoff File: new xlf name for b.md.md, keep the hyperlink
    # (same relationship id) but update its value + display text.
    $ws.Range("C3").Value = $loc.File
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$C$3') {
            $h.TextToDisplay = $loc.File
        }
    }

    # Latest Handoff Datetime: new handoff timestamp.
    $ws.Range("D3").Value = $loc.Datetime
}
